# Re-organise files to Model-View-Controller
# Append a new booking record (row 5) to the BookingHistory sheet, mirroring
# the existing rows (same customer/cinema/movie) with its own Quantity/Total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing row (A4:G4) down into row 5 so the new row
# inherits the same cell formatting/style and text-typed columns (e.g. the
# "Mobile" column stays text instead of being auto-coerced to a number).
$ws.Range("A4:G4").Copy($ws.Range("A5:G5"))

# Overwrite the per-booking figures for the new row.
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 48
